# Update the "取得日時" (acquired datetime) timestamp for all data rows
# on the active sheet (ランサーズ) from 2025-11-05 18:25:40 to 2025-11-05 18:33:29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A12").Value = "2025-11-05 18:33:29"
